$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename "user" -> "username" (A1). B1/C1 unchanged.
$ws.Range("A1").Value = "username"

# Row 2
$ws.Range("A2").Value = "testuser"
$ws.Range("B2").Value = "team@9%^"
$ws.Range("C2").Value = "Please Check the Username/Password"

# Row 3
$ws.Range("A3").Value = "Team9testuser"
$ws.Range("B3").Value = "team@9%^"
$ws.Range("C3").Value = "Please Check the Password"

# Row 4
$ws.Range("A4").Value = "testuser"
$ws.Range("B4").Value = "BddCoders@9"
$ws.Range("C4").Value = "Please Check the Username"

# Row 5
$ws.Range("B5").Value = "BddCoders@9"
$ws.Range("C5").Value = "Please Enter the Username"

# Row 6
$ws.Range("A6").Value = "Team9testuser"
$ws.Range("C6").Value = "Please Enter the Password"
$ws.Range("C6").Font.Size = 11

# Row 7
$ws.Range("C7").Value = "Please Enter the Username and Password"
$ws.Range("C7").Font.Size = 11

# Remove the single-cell selection that Excel previously remembered (B2)
$ws.Range("A1").Select() | Out-Null
